$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "Power Meter" description used by the two FooGadget rows.
$ws.Cells.Item(15, 6).Value = "Wireless power consumption meter"
$ws.Cells.Item(16, 6).Value = "Wireless power consumption meter"

# Row 24 - HueBridge
$ws.Cells.Item(24, 2).Value = "PHILIPS Hue"
$ws.Cells.Item(24, 3).Value = "RF Remote"
$ws.Cells.Item(24, 4).Value = "Zigbee"
$ws.Cells.Item(24, 6).Value = "Connects to the PHILPS Hue bridge device which controls a Hue system"

# Row 25 - HueLamp
$ws.Cells.Item(25, 2).Value = "PHILIPS Hue"
$ws.Cells.Item(25, 3).Value = "RF Remote"
$ws.Cells.Item(25, 4).Value = "Zigbee"
$ws.Cells.Item(25, 5).Value = "RF Bulb"
$ws.Cells.Item(25, 6).Value = "PHILPS Hue remote controlled lamp bulb"

# Row 26 - IkeaColorTemperatureLamp
$ws.Cells.Item(26, 2).Value = "IKEA Trådfri"
$ws.Cells.Item(26, 3).Value = "RF Remote"
$ws.Cells.Item(26, 4).Value = "Zigbee"
$ws.Cells.Item(26, 5).Value = "RF Bulb"
$ws.Cells.Item(26, 6).Value = "IKEA Trådfri remote controlled lamp bulb with adjustable color temperature"

# Row 27 - IkeaGateway
$ws.Cells.Item(27, 2).Value = "IKEA Trådfri"
$ws.Cells.Item(27, 3).Value = "RF Remote"
$ws.Cells.Item(27, 4).Value = "Zigbee"
$ws.Cells.Item(27, 6).Value = "IKEA Trådfri remote gateway"

# Row 28 - IkeaLamp
$ws.Cells.Item(28, 2).Value = "IKEA Trådfri"
$ws.Cells.Item(28, 3).Value = "RF Remote"
$ws.Cells.Item(28, 4).Value = "Zigbee"
$ws.Cells.Item(28, 5).Value = "RF Bulb"
$ws.Cells.Item(28, 6).Value = "IKEA Trådfri remote controlled lamp bulb"

# Row 29 - IntervalTimer
$ws.Cells.Item(29, 2).Value = "Timers"
$ws.Cells.Item(29, 3).Value = "Automation"
$ws.Cells.Item(29, 6).Value = "Delay timer with repeating actions"

# Row 30 - JeeLink
$ws.Cells.Item(30, 2).Value = "433MHz"
$ws.Cells.Item(30, 3).Value = "RF Remote"
$ws.Cells.Item(30, 4).Value = "JeeLink"
$ws.Cells.Item(30, 6).Value = "Item and firmware to use JeeLink Classic as RF Transmitter for 433MHz remote control"

# Move the selection/view like the author left it: cursor on B31, no frozen scroll offset.
$ws.Range("B31").Select()
